$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1: new header cell F1 = "expectedcount" (Text-formatted, like B1:E1) ----
$ws.Range("F1").Value = "expectedcount"
$ws.Range("F1").NumberFormat = "@"

# ---- Row 2: D2/E2 change from numbers to text "20170603"; new F2 = "Error" ----
$ws.Range("D2").Value = "20170603"
$ws.Range("E2").Value = "20170603"
$ws.Range("F2").Value = "Error"

# ---- Row 3: B3 -> "hyderabad"; D3 -> text "20170608"; E3 -> text "20170610"; new F3 = 2 ----
$ws.Range("B3").Value = "hyderabad"
$ws.Range("D3").Value = "20170608"
$ws.Range("E3").Value = "20170610"
$ws.Range("F3").Value = 2

# ---- Row 4: brand-new row ----
$ws.Range("A4").Value = "json"
$ws.Range("B4").Value = "mumbai"
$ws.Range("C4").Value = "goa"

# D4 must stay a genuine NUMBER (20170608) even though the column/cell uses the
# "Text" number format (numFmtId 49). Temporarily switch the cell to the Normal
# style (General format) so the numeric literal is not coerced to text, assign
# the value, then restore the Text format used by the rest of the column.
$ws.Range("D4").Style = "Normal"
$ws.Range("D4").Value = 20170608
$ws.Range("D4").NumberFormat = "@"

$ws.Range("E4").Value = "20170610"
$ws.Range("F4").Value = 2

# ---- Column widths ----
# Column B widens slightly to fit "hyderabad"; new column F gets a width too.
# (ColumnWidth is quantized internally, so we pick the input that lands on the
# closest achievable rendered width.)
$ws.Columns.Item(2).ColumnWidth = 9.63
$ws.Columns.Item(6).ColumnWidth = 13.6

# ---- Selection moves to F3 ----
$ws.Range("F3").Select()
